# end of sprint 2 commit
# Adds a new "Search Report Path" setting row to both the Assets and
# Local Config sheets, just above the existing "Retry Number" row, and
# grows each sheet's table to include the new row.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "Assets" ------------------------------------------------
$wsAssets = $wb.Worksheets.Item(1)
$tblAssets = $wsAssets.ListObjects.Item(1)

# Insert a new row above the last data row ("Retry Number" @ row 19),
# then grow the table range to cover it.
$wsAssets.Rows.Item(19).Insert()
$tblAssets.Resize($wsAssets.Range("A1:B20"))

$wsAssets.Range("A19").Value = "Search Report Path"
$wsAssets.Range("B19").Value = "Search Report Path"

# ---- Sheet 2: "Local Config" ------------------------------------------
$wsConfig = $wb.Worksheets.Item(2)
$tblConfig = $wsConfig.ListObjects.Item(1)

# Insert a new row above the last data row ("Retry Number" @ row 18),
# then grow the table range to cover it.
$wsConfig.Rows.Item(18).Insert()
$tblConfig.Resize($wsConfig.Range("A1:B19"))

$wsConfig.Range("A18").Value = "Search Report Path"
$wsConfig.Range("B18").Value = "C:\temp\Search Report.xlsx"
$wsConfig.Range("B18").WrapText = $true

# ---- Selection / active sheet ------------------------------------------
$wsAssets.Range("B19").Select() | Out-Null
$wsConfig.Activate() | Out-Null
$wsConfig.Range("E21").Select() | Out-Null

$wb.Save()
